$wb = $excel.ActiveWorkbook

# Sheet "Users to Delete" holds the Database ID (UUID) values in column D
$wsUsers = $wb.Worksheets.Item("Users to Delete")
$wsUsers.Range("D2").Value = "2e1afb38-45a5-41ae-9b63-9ab2c2c755cb"
$wsUsers.Range("D3").Value = "a760c524-bdb6-4454-9939-02e59dd7031d"
$wsUsers.Range("D4").Value = "051fc0a9-91fc-4865-84ed-d48d92026fb3"
$wsUsers.Range("D5").Value = "e522f464-167a-47f9-af6b-5de82126ad80"

# Sheet "Summary" holds the "Report Generated" timestamp in B6
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B6").Value = "11/10/2025, 3:18:26 PM"
